# Weekly data refresh: shift all forecast week-start dates forward by one
# week (rows 2-17 on "Forecast Comparison"), bump the Amazon P80 forecast
# for week 1, and refresh the derived date labels on "Summary".

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Forecast Comparison" ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# Force text formatting on the target cells first so the new date-looking
# strings are kept as plain text (matching how they were stored before),
# instead of being auto-converted into Excel date serial numbers.
$ws1.Range("B2:B17").NumberFormat = "@"

$ws1.Range("B2").Value = "2025-02-02"
$ws1.Range("B3").Value = "2025-02-09"
$ws1.Range("B4").Value = "2025-02-16"
$ws1.Range("B5").Value = "2025-02-23"
$ws1.Range("B6").Value = "2025-03-02"
$ws1.Range("B7").Value = "2025-03-09"
$ws1.Range("B8").Value = "2025-03-16"
$ws1.Range("B9").Value = "2025-03-23"
$ws1.Range("B10").Value = "2025-03-30"
$ws1.Range("B11").Value = "2025-04-06"
$ws1.Range("B12").Value = "2025-04-13"
$ws1.Range("B13").Value = "2025-04-20"
$ws1.Range("B14").Value = "2025-04-27"
$ws1.Range("B15").Value = "2025-05-04"
$ws1.Range("B16").Value = "2025-05-11"
$ws1.Range("B17").Value = "2025-05-18"

# Amazon P80 Forecast for week 1 (W1) increased from 1 to 2
$ws1.Range("G2").Value = 2

# --- Sheet 2: "Summary" ---
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B2").NumberFormat = "@"
$ws2.Range("B13").NumberFormat = "@"
$ws2.Range("B15").NumberFormat = "@"

$ws2.Range("B2").Value = "2023-01-01 to 2025-01-26"
$ws2.Range("B13").Value = "2025-02-02"
$ws2.Range("B15").Value = "2025-02-09"
